# Rename sheet "strategy_id-5008" to "strategy_id-5007", then add a new
# sheet "strategy_id-5009" (a duplicate of strategy_id-5007) right after it.

$wb = $excel.ActiveWorkbook

$target = $wb.Worksheets.Item("strategy_id-5008")
$target.Name = "strategy_id-5007"

# Duplicate the (renamed) sheet, placing the copy immediately after it.
$target.Copy($null, $target)

$newSheet = $wb.Worksheets.Item($target.Index + 1)
$newSheet.Name = "strategy_id-5009"
